$wb = $excel.ActiveWorkbook

$wsActie = $wb.Worksheets.Item("Actieplan")
$wsUren  = $wb.Worksheets.Item("Werkelijke uren")

# Add the two new logbook entries (hours) on the "Werkelijke uren" sheet.
$wsUren.Range("C15").Value = 1.75
$wsUren.Range("C16").Value = 3

# Update selection on each sheet to match the final state.
$wsActie.Range("G6").Select()
$wsUren.Range("C16").Select()

# Make "Werkelijke uren" the active (selected) sheet/tab.
$wsUren.Activate()
